$wb = $excel.ActiveWorkbook

# zh-cn sheet: update "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) for the row referencing
# 774b3b7d-e236-4366-a38c-c53fbc9656e2.f4f6b43204ea017f3f102a551408c4396ca2b361.zh-cn.xlf
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-25 08:47:40"
$wsZhCn.Range("H2").Value = "2016-03-25 08:48:20"

# de-de sheet: update "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) for the row referencing
# 774b3b7d-e236-4366-a38c-c53fbc9656e2.f4f6b43204ea017f3f102a551408c4396ca2b361.de-de.xlf
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-25 08:47:49"
$wsDeDe.Range("H2").Value = "2016-03-25 08:48:39"
